$d = $word.ActiveDocument

$searchRange = $d.Content
$found = $searchRange.Find.Execute("match with the captions ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
Write-Output "Found: $found"
if ($found) {
    $searchRange.Collapse(0)  # wdCollapseEnd
    $searchRange.InsertAfter("(Jackie)")
    $searchRange.Font.Color = 0
    $searchRange.Font.TextColor.ObjectThemeColor = 13
    Write-Output "Inserted"
}
